# "Stats: cli graphs and updated overview"
#
# The workbook opens with "Sheet1" (sheet2.xml on disk) as the active sheet
# (activeTab="1" / tabSelected="1"), so $wb.ActiveSheet already resolves to it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated overview numbers on the "Sheet1" stats table.
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = 2
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 3

# The failure marker in F5 changes from a generic "x" to "OOM".
$ws.Range("F5").Value = "OOM"

# Move the live selection to reflect where the author ended up (cli graphs).
$ws.Range("I24").Select() | Out-Null
